$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header A1 from "Issue Date" to "Date"
$ws.Range("A1").Value = "Date"

# Add new row 6 of data. The sheet stores every value as plain text
# (even numbers/dates), so force text format ("@") on the cells whose
# content would otherwise be auto-detected as a number or a date -
# this keeps them stored as literal text, matching the rest of the
# sheet.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("L6").NumberFormat = "@"

$ws.Range("A6").Value = "2023-03-15"
$ws.Range("B6").Value = "INV052"
$ws.Range("C6").Value = "Tech Enterprise"
$ws.Range("D6").Value = "Project X"
$ws.Range("E6").Value = "200"
$ws.Range("F6").Value = "30"
$ws.Range("G6").Value = "EUR"
$ws.Range("H6").Value = "6000"
$ws.Range("I6").Value = "1200"
$ws.Range("J6").Value = "7200"
$ws.Range("K6").Value = "Paid"
$ws.Range("L6").Value = "2023-03-20"
$ws.Range("M6").Value = "Sales"
$ws.Range("N6").Value = "Receipt"
$ws.Range("O6").Value = "Bank Transfer"
